# Generate Report for Handback
# Update the timestamp strings that reflect the latest handoff/handback
# xliff-generation times, as produced by a fresh run of the report generator.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 2df91cee...md (Overview row 2)
# and the matching "Correspond Handoff Datetime" for the de-de locale
# row 2 (same shared timestamp: the de-de xliff was generated last).
$wsOverview.Range("G2").Value = "2016-08-26 23:04:45"
$wsDeDe.Range("H2").Value     = "2016-08-26 23:04:45"

# zh-cn locale row 2: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-26 23:04:40"
$wsZhCn.Range("K2").Value = "2016-08-26 23:04:58"

# de-de locale row 2: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-26 23:05:15"
